$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: Dawuane Smoot / Group1
$ws.Range("A5:F5").Copy()
$ws.Range("A11:F11").PasteSpecial(-4122)
$ws.Range("A11").Value = "Dawuane Smoot"
$ws.Range("B11").Value = "Group1"
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 26
$ws.Range("E11").Value = 18.33333333333333
$ws.Range("F11").Value = 7.666666666666667

# Row 12: Dawuane Smoot / Group2
$ws.Range("A5:F5").Copy()
$ws.Range("A12:F12").PasteSpecial(-4122)
$ws.Range("A12").Value = "Dawuane Smoot"
$ws.Range("B12").Value = "Group2"
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = 15.66666666666667
$ws.Range("E12").Value = 11.33333333333333
$ws.Range("F12").Value = 4.333333333333333

# Row 13: Dawuane Smoot / Difference
$ws.Range("A5:F5").Copy()
$ws.Range("A13:F13").PasteSpecial(-4122)
$ws.Range("A13").Value = "Dawuane Smoot"
$ws.Range("B13").Value = "Difference"
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = -10.33333333333333
$ws.Range("E13").Value = -6.999999999999998
$ws.Range("F13").Value = -3.333333333333334

# Row 14: DeMarcus Walker / Group1
$ws.Range("A2:F2").Copy()
$ws.Range("A14:F14").PasteSpecial(-4122)
$ws.Range("A14").Value = "DeMarcus Walker"
$ws.Range("B14").Value = "Group1"
$ws.Range("C14").Value = 0.6666666666666666
$ws.Range("D14").Value = 23.66666666666667
$ws.Range("E14").Value = 14.33333333333333
$ws.Range("F14").Value = 9.333333333333334

# Row 15: DeMarcus Walker / Group2
$ws.Range("A2:F2").Copy()
$ws.Range("A15:F15").PasteSpecial(-4122)
$ws.Range("A15").Value = "DeMarcus Walker"
$ws.Range("B15").Value = "Group2"
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 36.33333333333334
$ws.Range("E15").Value = 21.66666666666667
$ws.Range("F15").Value = 14.66666666666667

# Row 16: DeMarcus Walker / Difference
$ws.Range("A2:F2").Copy()
$ws.Range("A16:F16").PasteSpecial(-4122)
$ws.Range("A16").Value = "DeMarcus Walker"
$ws.Range("B16").Value = "Difference"
$ws.Range("C16").Value = -0.6666666666666666
$ws.Range("D16").Value = 12.66666666666667
$ws.Range("E16").Value = 7.333333333333334
$ws.Range("F16").Value = 5.333333333333332

# Row 17: Jacob Martin / Group1
$ws.Range("A5:F5").Copy()
$ws.Range("A17:F17").PasteSpecial(-4122)
$ws.Range("A17").Value = "Jacob Martin"
$ws.Range("B17").Value = "Group1"
$ws.Range("C17").Value = 1.333333333333333
$ws.Range("D17").Value = 18
$ws.Range("E17").Value = 10
$ws.Range("F17").Value = 8

# Row 18: Jacob Martin / Group2
$ws.Range("A5:F5").Copy()
$ws.Range("A18:F18").PasteSpecial(-4122)
$ws.Range("A18").Value = "Jacob Martin"
$ws.Range("B18").Value = "Group2"
$ws.Range("C18").Value = 0.3333333333333333
$ws.Range("D18").Value = 10.44444444444444
$ws.Range("E18").Value = 6.222222222222222
$ws.Range("F18").Value = 4.222222222222222

# Row 19: Jacob Martin / Difference
$ws.Range("A5:F5").Copy()
$ws.Range("A19:F19").PasteSpecial(-4122)
$ws.Range("A19").Value = "Jacob Martin"
$ws.Range("B19").Value = "Difference"
$ws.Range("C19").Value = -1
$ws.Range("D19").Value = -7.555555555555555
$ws.Range("E19").Value = -3.777777777777778
$ws.Range("F19").Value = -3.777777777777778

# Row 20: Kenny Clark / Group1
$ws.Range("A2:F2").Copy()
$ws.Range("A20:F20").PasteSpecial(-4122)
$ws.Range("A20").Value = "Kenny Clark"
$ws.Range("B20").Value = "Group1"
$ws.Range("C20").Value = 0.3333333333333333
$ws.Range("D20").Value = 50.66666666666666
$ws.Range("E20").Value = 28.66666666666667
$ws.Range("F20").Value = 22

# Row 21: Kenny Clark / Group2
$ws.Range("A2:F2").Copy()
$ws.Range("A21:F21").PasteSpecial(-4122)
$ws.Range("A21").Value = "Kenny Clark"
$ws.Range("B21").Value = "Group2"
$ws.Range("C21").Value = 1.666666666666667
$ws.Range("D21").Value = 44.66666666666666
$ws.Range("E21").Value = 23
$ws.Range("F21").Value = 21.66666666666667

# Row 22: Kenny Clark / Difference
$ws.Range("A2:F2").Copy()
$ws.Range("A22:F22").PasteSpecial(-4122)
$ws.Range("A22").Value = "Kenny Clark"
$ws.Range("B22").Value = "Difference"
$ws.Range("C22").Value = 1.333333333333333
$ws.Range("D22").Value = -6
$ws.Range("E22").Value = -5.666666666666668
$ws.Range("F22").Value = -0.3333333333333321

# Row 23: Maliek Collins / Group1
$ws.Range("A5:F5").Copy()
$ws.Range("A23:F23").PasteSpecial(-4122)
$ws.Range("A23").Value = "Maliek Collins"
$ws.Range("B23").Value = "Group1"
$ws.Range("C23").Value = 0.3333333333333333
$ws.Range("D23").Value = 21.33333333333333
$ws.Range("E23").Value = 13.66666666666667
$ws.Range("F23").Value = 7.666666666666667

# Row 24: Maliek Collins / Group2
$ws.Range("A5:F5").Copy()
$ws.Range("A24:F24").PasteSpecial(-4122)
$ws.Range("A24").Value = "Maliek Collins"
$ws.Range("B24").Value = "Group2"
$ws.Range("C24").Value = 1
$ws.Range("D24").Value = 37
$ws.Range("E24").Value = 22
$ws.Range("F24").Value = 15

# Row 25: Maliek Collins / Difference
$ws.Range("A5:F5").Copy()
$ws.Range("A25:F25").PasteSpecial(-4122)
$ws.Range("A25").Value = "Maliek Collins"
$ws.Range("B25").Value = "Difference"
$ws.Range("C25").Value = 0.6666666666666667
$ws.Range("D25").Value = 15.66666666666667
$ws.Range("E25").Value = 8.333333333333334
$ws.Range("F25").Value = 7.333333333333333

$excel.CutCopyMode = 0
